# Fruta / hortaliza, semanal
# The weekly refresh reshuffled which data row holds which record. Columns
# A,B,C,E,F,G,H,I,J are identical for every row in this subconjunto (same
# market/product), so only D (Fecha), K..T (Variedad..Kg/unidad) actually
# move between rows. Capture the "before" snapshot for those columns, then
# redistribute it across rows 2-16 according to the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row -> old row it now contains the data of
$mapping = @{
    2  = 13
    3  = 8
    4  = 9
    5  = 2
    6  = 15
    7  = 5
    8  = 6
    9  = 16
    10 = 3
    11 = 4
    12 = 7
    13 = 14
    14 = 12
    15 = 10
    16 = 11
}

$cols = @(4, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)  # D, K, L, M, N, O, P, Q, R, S, T

# Snapshot current ("before") values for the moving columns, rows 2-16
$snapshot = @{}
for ($r = 2; $r -le 16; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write new values based on the permutation mapping
for ($r = 2; $r -le 16; $r++) {
    $src = $mapping[$r]
    $rowVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c]
    }
}
